$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh rows 2-7 with updated TPM-derived NATMI values (FAPs/MuSCs as sending clusters for Alcam-Nrp1)
# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Alcam"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5683613333333334
$ws.Range("H2").Value = 1.705084
$ws.Range("I2").Value = 0.4361027177196302
$ws.Range("J2").Value = 0.4361027177196302
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 62.88503704249067
$ws.Range("R2").Value = 565.9653333824159
$ws.Range("S2").Value = 0.2388281176732297
$ws.Range("T2").Value = 0.2388281176732297

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Alcam"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5683613333333334
$ws.Range("H3").Value = 1.705084
$ws.Range("I3").Value = 0.4361027177196302
$ws.Range("J3").Value = 0.4361027177196302
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 36.18337632303778
$ws.Range("R3").Value = 325.65038690734
$ws.Range("S3").Value = 0.1374191391897279
$ws.Range("T3").Value = 0.1374191391897279

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Alcam"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5683613333333334
$ws.Range("H4").Value = 1.705084
$ws.Range("I4").Value = 0.4361027177196302
$ws.Range("J4").Value = 0.4361027177196302
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.188324
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 15.76034224880178
$ws.Range("R4").Value = 141.843080239216
$ws.Range("S4").Value = 0.05985546085667257
$ws.Range("T4").Value = 0.05985546085667257

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Alcam"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7349126666666667
$ws.Range("H5").Value = 2.204738
$ws.Range("I5").Value = 0.5638972822803697
$ws.Range("J5").Value = 0.5638972822803697
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 81.31272758350133
$ws.Range("R5").Value = 731.8145482515118
$ws.Range("S5").Value = 0.3088137748654266
$ws.Range("T5").Value = 0.3088137748654266

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Alcam"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7349126666666667
$ws.Range("H6").Value = 2.204738
$ws.Range("I6").Value = 0.5638972822803697
$ws.Range("J6").Value = 0.5638972822803697
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 46.78647195545889
$ws.Range("R6").Value = 421.07824759913
$ws.Range("S6").Value = 0.1776881362436586
$ws.Range("T6").Value = 0.1776881362436586

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Alcam"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7349126666666667
$ws.Range("H7").Value = 2.204738
$ws.Range("I7").Value = 0.5638972822803697
$ws.Range("J7").Value = 0.5638972822803697
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.188324
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 20.37871767545689
$ws.Range("R7").Value = 183.408459079112
$ws.Range("S7").Value = 0.07739537117128455
$ws.Range("T7").Value = 0.07739537117128455

# Drop the former MuSCs sending-cluster block (old rows 8-10); its data has been merged into rows 2-7 above
$ws.Rows("8:10").Delete()
